$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 554.5640449799788
$ws.Range("D2").Value = 136.2932116927976
$ws.Range("G2").Value = 514
$ws.Range("C3").Value = 37.64364162057874
$ws.Range("D3").Value = 6.502298985304678
$ws.Range("F3").Value = 32.95
$ws.Range("G3").Value = 37.95
$ws.Range("H3").Value = 41.59
$ws.Range("C4").Value = 2.002309828873187
$ws.Range("D4").Value = 2.559219627636518
$ws.Range("F4").Value = 0.63
$ws.Range("G4").Value = 1.29
$ws.Range("C5").Value = 323.2802032006397
$ws.Range("D5").Value = 10.88097831185358
$ws.Range("F5").Value = 316.25
$ws.Range("G5").Value = 324.79
$ws.Range("H5").Value = 331.56
$ws.Range("C6").Value = 21.18694478796489
$ws.Range("D6").Value = 2.589795419055544
$ws.Range("F6").Value = 19.71
$ws.Range("G6").Value = 21.13
$ws.Range("H6").Value = 22.52
$ws.Range("C7").Value = -76.8894234526639
$ws.Range("D7").Value = 22.93711828158775
$ws.Range("C8").Value = 7.415725274312365
$ws.Range("D8").Value = 7.102650580236241
$ws.Range("C9").Value = 9.322666065952859
$ws.Range("D9").Value = 1.685729385373283
$ws.Range("C10").Value = 867.83009024499
$ws.Range("D10").Value = 0.4614888068895844
$ws.Range("C11").Value = 0.555930376485119
$ws.Range("D11").Value = 0.5890373160311491
$ws.Range("C12").Value = 22.74735387100284
$ws.Range("D12").Value = 12.29392615808524
$ws.Range("C13").Value = 0.6740300081065715
$ws.Range("D13").Value = 0.7506847753406487
$ws.Range("C14").Value = 1.827769834665784
$ws.Range("D14").Value = 1.66448662828912
$ws.Range("C15").Value = 94.14942345266408
$ws.Range("D15").Value = 22.93711828158775
$ws.Range("C16").Value = -86.00904533313746
$ws.Range("D16").Value = 20.50881087638902
$ws.Range("F16").Value = -102.7900974965257
$ws.Range("G16").Value = -84.26572375596102
$ws.Range("H16").Value = -70.33779541063677
$ws.Range("C17").Value = -78.59332005882511
$ws.Range("D17").Value = 25.48295048108626
$ws.Range("F17").Value = -93.93380807687734
$ws.Range("G17").Value = -73.71081852649533
$ws.Range("H17").Value = -60.33195619988427
